$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4033.7273
$ws.Range("I76").Value = 3430.889
$ws.Range("K76").Value = 3430.889
$ws.Range("M76").Value = -3115.889
$ws.Range("H79").Value = 4033.7273
$ws.Range("I79").Value = 3430.889
$ws.Range("K79").Value = 3430.889
$ws.Range("M79").Value = -2338.889
$ws.Range("H86").Value = 2873.3333
$ws.Range("I86").Value = 2865.7144
$ws.Range("K86").Value = 2865.7144
$ws.Range("M86").Value = -1742.7144
$ws.Range("H88").Value = 1771.4375
$ws.Range("J88").Value = 1735.4286
$ws.Range("L88").Value = 1735.4286
$ws.Range("N88").Value = -2547.4286
$ws.Range("H89").Value = 2873.3333
$ws.Range("I89").Value = 2865.7144
$ws.Range("K89").Value = 14328.572
$ws.Range("M89").Value = -8712.572
$ws.Range("H91").Value = 1771.4375
$ws.Range("J91").Value = 1735.4286
$ws.Range("L91").Value = 1735.4286
$ws.Range("N91").Value = -4543.4286
$ws.Range("H116").Value = 5416.6665
$ws.Range("J116").Value = 5000
$ws.Range("L116").Value = 5000
$ws.Range("N116").Value = -11884
$ws.Range("H135").Value = 34412.2
$ws.Range("I135").Value = 980.7619
$ws.Range("K135").Value = 8826.857099999999
$ws.Range("M135").Value = -6291.857099999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 48738.754
$ws.Range("I32").Value = 28316.691
$ws.Range("J32").Value = 181482.17
$ws.Range("K32").Value = 28316.691
$ws.Range("L32").Value = 181482.17
$ws.Range("M32").Value = -28029.691
$ws.Range("N32").Value = -182056.17
$ws.Range("H74").Value = 1318.6531
$ws.Range("I74").Value = 597.05554
$ws.Range("J74").Value = 3316.923
$ws.Range("K74").Value = 597.05554
$ws.Range("L74").Value = 3316.923
$ws.Range("M74").Value = 276.94446
$ws.Range("N74").Value = -5064.923
$ws.Range("H77").Value = 1318.6531
$ws.Range("I77").Value = 597.05554
$ws.Range("J77").Value = 3316.923
$ws.Range("K77").Value = 2985.2777
$ws.Range("L77").Value = 16584.615
$ws.Range("M77").Value = 1382.7223
$ws.Range("N77").Value = -25320.615
$ws.Range("H122").Value = 2068.1667
$ws.Range("I122").Value = 2068.1667
$ws.Range("K122").Value = 6204.500100000001
$ws.Range("M122").Value = -3754.500100000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3375.2
$ws.Range("I86").Value = 3092
$ws.Range("K86").Value = 3092
$ws.Range("M86").Value = -1969
$ws.Range("H89").Value = 3375.2
$ws.Range("I89").Value = 3092
$ws.Range("K89").Value = 15460
$ws.Range("M89").Value = -9844
$ws.Range("H96").Value = 30214
$ws.Range("I96").Value = 10428
$ws.Range("J96").Value = 50000
$ws.Range("K96").Value = 10428
$ws.Range("L96").Value = 50000
$ws.Range("M96").Value = -7682
$ws.Range("N96").Value = -55492

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 889
$ws.Range("I58").Value = 695
$ws.Range("K58").Value = 695
$ws.Range("M58").Value = -492
$ws.Range("H62").Value = 12317.454
$ws.Range("J62").Value = 2932.6667
$ws.Range("L62").Value = 2932.6667
$ws.Range("N62").Value = -4180.6667
$ws.Range("H65").Value = 12317.454
$ws.Range("J65").Value = 2932.6667
$ws.Range("L65").Value = 14663.3335
$ws.Range("N65").Value = -20903.3335
$ws.Range("H136").Value = 889
$ws.Range("I136").Value = 695
$ws.Range("K136").Value = 2085
$ws.Range("M136").Value = 465

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 541.1667
$ws.Range("J113").Value = 499.4
$ws.Range("L113").Value = 1498.2
$ws.Range("N113").Value = -5838.2
$ws.Range("H132").Value = 1042.3529
$ws.Range("I132").Value = 891.9091
$ws.Range("J132").Value = 1318.1666
$ws.Range("K132").Value = 8027.1819
$ws.Range("L132").Value = 11863.4994
$ws.Range("M132").Value = -5497.1819
$ws.Range("N132").Value = -16923.4994
$ws.Range("H134").Value = 1664.1111
$ws.Range("I134").Value = 1664.1111
$ws.Range("K134").Value = 4992.3333
$ws.Range("M134").Value = 77.66669999999976

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5275.6665
$ws.Range("J80").Value = 6997.75
$ws.Range("L80").Value = 6997.75
$ws.Range("N80").Value = -8993.75
$ws.Range("H83").Value = 5275.6665
$ws.Range("J83").Value = 6997.75
$ws.Range("L83").Value = 34988.75
$ws.Range("N83").Value = -44972.75
$ws.Range("H132").Value = 2501.3333
$ws.Range("I132").Value = 2541.48
$ws.Range("J132").Value = 1999.5
$ws.Range("K132").Value = 7624.440000000001
$ws.Range("L132").Value = 5998.5
$ws.Range("M132").Value = -5094.440000000001
$ws.Range("N132").Value = -11058.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 40183.184
$ws.Range("I46").Value = 40183.184
$ws.Range("K46").Value = 40183.184
$ws.Range("M46").Value = -39995.184
$ws.Range("I68").Value = 1831.4
$ws.Range("J68").Value = 2796
$ws.Range("K68").Value = 1831.4
$ws.Range("L68").Value = 2796
$ws.Range("M68").Value = -1082.4
$ws.Range("N68").Value = -4294
$ws.Range("I71").Value = 1831.4
$ws.Range("J71").Value = 2796
$ws.Range("K71").Value = 9157
$ws.Range("L71").Value = 13980
$ws.Range("M71").Value = -5413
$ws.Range("N71").Value = -21468
$ws.Range("H74").Value = 23173.2
$ws.Range("I74").Value = 15000
$ws.Range("J74").Value = 25216.5
$ws.Range("K74").Value = 15000
$ws.Range("L74").Value = 25216.5
$ws.Range("M74").Value = -14002
$ws.Range("N74").Value = -27212.5
$ws.Range("H77").Value = 23173.2
$ws.Range("I77").Value = 15000
$ws.Range("J77").Value = 25216.5
$ws.Range("K77").Value = 45000
$ws.Range("L77").Value = 75649.5
$ws.Range("M77").Value = -40008
$ws.Range("N77").Value = -85633.5
$ws.Range("H82").Value = 3873.6667
$ws.Range("J82").Value = 3873.6667
$ws.Range("L82").Value = 3873.6667
$ws.Range("N82").Value = -4595.6667
$ws.Range("H85").Value = 3873.6667
$ws.Range("J85").Value = 3873.6667
$ws.Range("L85").Value = 3873.6667
$ws.Range("N85").Value = -6369.6667

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 26984
$ws.Range("I99").Value = 20000
$ws.Range("K99").Value = 20000
$ws.Range("M99").Value = -7682
